$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: insert a new "16-dec" column before column ES ---
# This shifts the existing ES:FW block (01-oct. ... 31-oct. and all the
# underlying daily price data for rows 2-25) one column to the right
# (ET:FX), matching the canonical diff.
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("ES").Insert()

# New header cell for the inserted column.
$wsSpot.Range("ES1").Value = "16-dec"

# The inserted column has no data yet for this new date, so every data row
# (2 through 25) gets the same placeholder used elsewhere in the sheet for
# missing values.
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 149).Value = "-"
}

# --- "Gaz" sheet: update the two most recent price points ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B178").Value = 25.93
$wsGaz.Range("B179").Value = 25.93
